$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM/CPL sheet gained two new placement rows, "C15" and "C16",
# inserted right after the existing "C14" row (old row 15 / new row 15),
# pushing every subsequent designator row (old D1..U4, old rows 16-78)
# down by two (new rows 18-80). The two trailing blank spacer rows
# (old rows 79-80) move down to new rows 81-82.

# Insert two blank rows at row 16, shifting row 16 and everything below
# it down by two rows. This preserves per-cell styles/number formats
# because Excel carries the row-16 formatting into the newly inserted
# rows.
$ws.Rows("16:17").Insert()

# Fill in the two new rows with the new designator data.
$ws.Range("A16").Value = "C15"
$ws.Range("B16").Value = -6.75
$ws.Range("C16").Value = 87.75
$ws.Range("D16").Value = "top"
$ws.Range("E16").Value = 270.0

$ws.Range("A17").Value = "C16"
$ws.Range("B17").Value = -7.0
$ws.Range("C17").Value = 28.5
$ws.Range("D17").Value = "top"
$ws.Range("E17").Value = 270.0

# The two trailing placeholder rows used to be rows 79-80 (each with just
# an "E" cell carrying the blank/border style) and are now rows 81-82.
# They picked up a matching "D" placeholder cell with the same style as
# their "E" cell. Copy the style from E81/E82 onto D81/D82 (format-only
# paste so we don't disturb the shared style table with a duplicate).
$ws.Range("E81").Copy() | Out-Null
$ws.Range("D81").PasteSpecial(-4122) | Out-Null

$ws.Range("E82").Copy() | Out-Null
$ws.Range("D82").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
